# Add logging and regenerate Excel with correct schema (RTM column)
Write-Host "Starting auction data schema update..."

$wb = $excel.ActiveWorkbook

# --- Sheet 1: "Available Players" -----------------------------------------
Write-Host "Updating sheet 'Available Players': adding Franchise ID column"
$ws1 = $wb.Worksheets.Item("Available Players")

# New header for column E
$ws1.Cells.Item(1, 5).Value = "Franchise ID"

# Per-player Franchise ID values (row 2..21)
$ws1.Cells.Item(2, 5).Value = 3
$ws1.Cells.Item(3, 5).Value = 1
$ws1.Cells.Item(4, 5).Value = 1
$ws1.Cells.Item(5, 5).Value = 8
$ws1.Cells.Item(6, 5).Value = 6
$ws1.Cells.Item(7, 5).Value = 8
$ws1.Cells.Item(8, 5).Value = 1
$ws1.Cells.Item(9, 5).Value = 2
$ws1.Cells.Item(10, 5).Value = 5
$ws1.Cells.Item(11, 5).Value = 7
$ws1.Cells.Item(12, 5).Value = 6
$ws1.Cells.Item(13, 5).Value = 8
$ws1.Cells.Item(14, 5).Value = 3
$ws1.Cells.Item(15, 5).Value = 3
$ws1.Cells.Item(16, 5).Value = 5
$ws1.Cells.Item(17, 5).Value = 7
$ws1.Cells.Item(18, 5).Value = 4
$ws1.Cells.Item(19, 5).Value = 7
$ws1.Cells.Item(20, 5).Value = 4
$ws1.Cells.Item(21, 5).Value = 1

# Match column D's width (15 characters) for the new column E
$ws1.Columns.Item(5).ColumnWidth = 14.14

Write-Host "Sheet 'Available Players' updated: 20 rows x 5 columns"

# --- Sheet 2: "Sold Players" -----------------------------------------------
Write-Host "Updating sheet 'Sold Players': adding RTM Used column"
$ws2 = $wb.Worksheets.Item("Sold Players")

$ws2.Cells.Item(1, 7).Value = "RTM Used"

# Match columns 1/4's width (10 characters) for the new column G
$ws2.Columns.Item(7).ColumnWidth = 9.14

Write-Host "Sheet 'Sold Players' updated: header row now has 7 columns"

Write-Host "Auction data schema update complete."
